$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.0
$ws.Range("B3").Value = 4.0
$ws.Range("B4").Value = 3.0
$ws.Range("B5").Value = 0.0
$ws.Range("B6").Value = 0.0
